$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.193.40"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.862.96"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Formula = "'1.002"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Formula = "'0.7097"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Formula = "'240.53"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Formula = "'0.3074"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Formula = "'0.07687"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").Formula = "'24.92"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "1.848.30"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").Formula = "'0.7161"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Formula = "'5.210"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Formula = "'90.19"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "29.191.67"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Formula = "'5.850"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Formula = "'243.16"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Formula = "'0.000007790"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "2.111.21"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Formula = "'7.952"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").Formula = "'1.002"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Formula = "'0.1573"
$ws.Range("E25").Value = "  +6.19%  "
$ws.Range("D26").Formula = "'162.43"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Formula = "'8.896"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Formula = "'18.18"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Formula = "'1.324"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Formula = "'4.348"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Formula = "'4.092"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Formula = "'0.05183"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").Formula = "'1.900"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").Formula = "'1.174"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("D36").Formula = "'0.7272"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("D37").Formula = "'2.686"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").Formula = "'0.01845"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Formula = "'2.688"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "1.141.44"
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").Formula = "'0.8984"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("D42").Formula = "'6.083"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Formula = "'72.12"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Formula = "'101.45"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Formula = "'0.5269"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.006.06"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Formula = "'1.764"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").Formula = "'9.287"
$ws.Range("D51").Formula = "'2.859"
$ws.Range("E51").Value = "  -1.32%  "
